$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.88
$ws.Range("I2").Value = 2.63
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 3.5
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.38
$ws.Range("Q2").Value = 2.7
$ws.Range("R2").Value = 1.44
$ws.Range("Y2").Value = 12
$ws.Range("AC2").Value = 6
$ws.Range("AF2").Value = 81
$ws.Range("AH2").Value = 11
$ws.Range("AJ2").Value = 26
$ws.Range("AN2").Value = 4.75

# Row 3 updates
$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 3.8
$ws.Range("L3").Value = 4.33
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("X3").Value = 9
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 15
$ws.Range("AG3").Value = 9.5
$ws.Range("AK3").Value = 34
$ws.Range("AY3").Value = 34

# Row 5 updates
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.3
$ws.Range("I5").Value = 5.75
$ws.Range("K5").Value = 2.05
$ws.Range("S5").Value = 1.5
$ws.Range("T5").Value = 2.5
$ws.Range("Z5").Value = 12
$ws.Range("AD5").Value = 6.5
$ws.Range("AG5").Value = 11
$ws.Range("AH5").Value = 26
$ws.Range("AI5").Value = 19
$ws.Range("AJ5").Value = 67
$ws.Range("AK5").Value = 51
$ws.Range("AT5").Value = 2.5
$ws.Range("AW5").Value = 7
$ws.Range("AX5").Value = 34

# Row 6 updates
$ws.Range("G6").Value = 3.2
$ws.Range("H6").Value = 2.88
$ws.Range("I6").Value = 2.5
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("AE6").Value = 19
$ws.Range("AN6").Value = 4.75
$ws.Range("AO6").Value = 19
$ws.Range("AR6").Value = 101

# Row 7 updates
$ws.Range("G7").Value = 2.35
$ws.Range("I7").Value = 2.9
$ws.Range("J7").Value = 3.1
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.83
$ws.Range("X7").Value = 11
$ws.Range("Y7").Value = 10
$ws.Range("Z7").Value = 23
$ws.Range("AD7").Value = 6
$ws.Range("AE7").Value = 15
$ws.Range("AI7").Value = 11
$ws.Range("AP7").Value = 26
$ws.Range("AX7").Value = 17

# Row 8 updates
$ws.Range("G8").Value = 1.8
$ws.Range("I8").Value = 4.5
$ws.Range("L8").Value = 5
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("W8").Value = 6
$ws.Range("X8").Value = 8
$ws.Range("Z8").Value = 15
$ws.Range("AC8").Value = 8
$ws.Range("AE8").Value = 17
$ws.Range("AF8").Value = 51
$ws.Range("AL8").Value = 41
$ws.Range("AN8").Value = 3.75
$ws.Range("AR8").Value = 51
$ws.Range("AT8").Value = 2.63
$ws.Range("AU8").Value = 9
$ws.Range("AW8").Value = 6
$ws.Range("AX8").Value = 26
$ws.Range("AY8").Value = 34
$ws.Range("AZ8").Value = 81
$ws.Range("BA8").Value = 126
$ws.Range("BB8").Value = 301
